$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 309358.97
$ws.Range("I15").Value = 309358.97
$ws.Range("K15").Value = 928076.9099999999
$ws.Range("M15").Value = -927907.9099999999
$ws.Range("H17").Value = 1826.0769
$ws.Range("J17").Value = 1849.7368
$ws.Range("L17").Value = 5549.2104
$ws.Range("N17").Value = -5885.2104
$ws.Range("H40").Value = 4159.2
$ws.Range("I40").Value = 3836.2727
$ws.Range("K40").Value = 3836.2727
$ws.Range("M40").Value = -3661.2727
$ws.Range("H64").Value = 5304.737
$ws.Range("I64").Value = 5356.5
$ws.Range("J64").Value = 5159.8
$ws.Range("K64").Value = 5356.5
$ws.Range("L64").Value = 5159.8
$ws.Range("M64").Value = -5108.5
$ws.Range("N64").Value = -5655.8
$ws.Range("H67").Value = 5304.737
$ws.Range("I67").Value = 5356.5
$ws.Range("J67").Value = 5159.8
$ws.Range("K67").Value = 5356.5
$ws.Range("L67").Value = 5159.8
$ws.Range("M67").Value = -4498.5
$ws.Range("N67").Value = -6875.8
$ws.Range("H74").Value = 6421.0713
$ws.Range("I74").Value = 5809.9
$ws.Range("J74").Value = 7949
$ws.Range("K74").Value = 5809.9
$ws.Range("L74").Value = 7949
$ws.Range("M74").Value = -4873.9
$ws.Range("N74").Value = -9821
$ws.Range("H77").Value = 6421.0713
$ws.Range("I77").Value = 5809.9
$ws.Range("J77").Value = 7949
$ws.Range("K77").Value = 29049.5
$ws.Range("L77").Value = 39745
$ws.Range("M77").Value = -24369.5
$ws.Range("N77").Value = -49105
$ws.Range("H98").Value = 7937
$ws.Range("I98").Value = 7944.6665
$ws.Range("K98").Value = 7944.6665
$ws.Range("M98").Value = -6446.6665
$ws.Range("H122").Value = 7937
$ws.Range("I122").Value = 7944.6665
$ws.Range("K122").Value = 23833.9995
$ws.Range("M122").Value = -21383.9995
$ws.Range("H137").Value = 13025.046
$ws.Range("I137").Value = 4570.9165
$ws.Range("K137").Value = 13712.7495
$ws.Range("M137").Value = -11162.7495
$ws.Range("H138").Value = 6438.7666
$ws.Range("I138").Value = 3333.3333
$ws.Range("J138").Value = 6783.815
$ws.Range("K138").Value = 9999.999899999999
$ws.Range("L138").Value = 20351.445
$ws.Range("M138").Value = -4859.999899999999
$ws.Range("N138").Value = -30631.445
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3223.6667
$ws.Range("I61").Value = 3064.125
$ws.Range("K61").Value = 3064.125
$ws.Range("M61").Value = -2852.125
$ws.Range("H132").Value = 41669120
$ws.Range("I132").Value = 2585
$ws.Range("K132").Value = 7755
$ws.Range("M132").Value = -5225
$ws.Range("H136").Value = 3223.6667
$ws.Range("I136").Value = 3064.125
$ws.Range("K136").Value = 9192.375
$ws.Range("M136").Value = -6642.375
$ws.Range("H137").Value = 63324.082
$ws.Range("J137").Value = 79972.25
$ws.Range("L137").Value = 79972.25
$ws.Range("N137").Value = -90172.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3347.8572
$ws.Range("I20").Value = 3336.5
$ws.Range("J20").Value = 3363
$ws.Range("K20").Value = 3336.5
$ws.Range("L20").Value = 3363
$ws.Range("M20").Value = -3089.5
$ws.Range("N20").Value = -3857
$ws.Range("H82").Value = 9947.5
$ws.Range("J82").Value = 24283
$ws.Range("L82").Value = 24283
$ws.Range("N82").Value = -25049
$ws.Range("H85").Value = 9947.5
$ws.Range("J85").Value = 24283
$ws.Range("L85").Value = 24283
$ws.Range("N85").Value = -26935
$ws.Range("H86").Value = 4990.3
$ws.Range("J86").Value = 4640
$ws.Range("L86").Value = 4640
$ws.Range("N86").Value = -6886
$ws.Range("H89").Value = 4990.3
$ws.Range("J89").Value = 4640
$ws.Range("L89").Value = 23200
$ws.Range("N89").Value = -34432
$ws.Range("H94").Value = 129999.5
$ws.Range("J94").Value = 149999.5
$ws.Range("L94").Value = 149999.5
$ws.Range("N94").Value = -150901.5
$ws.Range("H132").Value = 91137.42999999999
$ws.Range("J132").Value = 91137.42999999999
$ws.Range("L132").Value = 91137.42999999999
$ws.Range("N132").Value = -101257.43
$ws.Range("H134").Value = 48621010
$ws.Range("I134").Value = 25010424
$ws.Range("J134").Value = 166673920
$ws.Range("K134").Value = 75031272
$ws.Range("L134").Value = 500021760
$ws.Range("M134").Value = -75028737
$ws.Range("N134").Value = -500026830
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6496084.5
$ws.Range("I31").Value = 6995245
$ws.Range("K31").Value = 6995245
$ws.Range("M31").Value = -6994950
$ws.Range("H34").Value = 6496084.5
$ws.Range("I34").Value = 6995245
$ws.Range("K34").Value = 6995245
$ws.Range("M34").Value = -6995043
$ws.Range("H74").Value = 47496.75
$ws.Range("J74").Value = 59993.5
$ws.Range("L74").Value = 59993.5
$ws.Range("N74").Value = -61741.5
$ws.Range("H77").Value = 47496.75
$ws.Range("J77").Value = 59993.5
$ws.Range("L77").Value = 179980.5
$ws.Range("N77").Value = -188716.5
$ws.Range("H132").Value = 5159.6924
$ws.Range("I132").Value = 4643.364
$ws.Range("K132").Value = 13930.092
$ws.Range("M132").Value = -11400.092
$ws.Range("H134").Value = 2779252
$ws.Range("I134").Value = 1509.2413
$ws.Range("J134").Value = 14287044
$ws.Range("K134").Value = 4527.7239
$ws.Range("L134").Value = 42861132
$ws.Range("M134").Value = -1992.7239
$ws.Range("N134").Value = -42866202
$ws.Range("H141").Value = 303476.66
$ws.Range("J141").Value = 350007.84
$ws.Range("L141").Value = 350007.84
$ws.Range("N141").Value = -360367.84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H113").Value = 1307.2858
$ws.Range("J113").Value = 1307.2858
$ws.Range("L113").Value = 3921.8574
$ws.Range("N113").Value = -8261.857400000001
$ws.Range("H128").Value = 154655.64
$ws.Range("I128").Value = 154655.64
$ws.Range("K128").Value = 463966.92
$ws.Range("M128").Value = -458986.92
$ws.Range("H137").Value = 14411
$ws.Range("J137").Value = 16688.666
$ws.Range("L137").Value = 50065.99800000001
$ws.Range("N137").Value = -60265.99800000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 32499.5
$ws.Range("J100").Value = 32499.5
$ws.Range("L100").Value = 32499.5
$ws.Range("N100").Value = -34663.5
$ws.Range("H122").Value = 2902.697
$ws.Range("I122").Value = 3210.2778
$ws.Range("J122").Value = 2533.6
$ws.Range("K122").Value = 9630.8334
$ws.Range("L122").Value = 7600.799999999999
$ws.Range("M122").Value = -7180.8334
$ws.Range("N122").Value = -12500.8
$ws.Range("H132").Value = 4021.375
$ws.Range("I132").Value = 5753
$ws.Range("J132").Value = 2289.75
$ws.Range("K132").Value = 17259
$ws.Range("L132").Value = 6869.25
$ws.Range("M132").Value = -14729
$ws.Range("N132").Value = -11929.25
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5210.5293
$ws.Range("I132").Value = 5361.25
$ws.Range("J132").Value = 2799
$ws.Range("K132").Value = 16083.75
$ws.Range("L132").Value = 8397
$ws.Range("M132").Value = -13553.75
$ws.Range("N132").Value = -13457
$ws.Range("H136").Value = 1401.0975
$ws.Range("I136").Value = 1067.1923
$ws.Range("J136").Value = 1979.8667
$ws.Range("K136").Value = 3201.5769
$ws.Range("L136").Value = 5939.6001
$ws.Range("M136").Value = -651.5769
$ws.Range("N136").Value = -11039.6001
$ws.Range("H140").Value = 69995
$ws.Range("I140").Value = 69995
$ws.Range("K140").Value = 69995
$ws.Range("M140").Value = -64815

Write-Output "Applied 200 cell updates across 7 sheets"